$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each parent currently has two email columns ("Email 1" / "Email 2"). Drop the
# second email column for each parent (originally F, J, N) and keep the first,
# renaming its header from "Parent N Email 1" to simply "Parent N Email".
# Delete right-to-left so the earlier column letters remain valid as we go.
$ws.Columns("N").Delete()
$ws.Columns("J").Delete()
$ws.Columns("F").Delete()

$ws.Range("E1").Value = "Parent 1 Email"
$ws.Range("H1").Value = "Parent 2 Email"
$ws.Range("K1").Value = "Parent 3 Email"

# Column deletion doesn't shift the worksheet's hyperlink anchors, so rebuild
# the hyperlink list to match the new layout: the two "Email 2" links
# (on the now-deleted columns) are dropped, and the remaining ones point at
# their new, shifted-left cells.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:mwenda.lilian@yahoo.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:mwenda.lilian@yahoo.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:Vusimuzi@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K3"), "mailto:info@bahari.org") | Out-Null

# Re-adding hyperlinks resets these cells' style; restore the original
# "Hyperlink" look (underlined, themed link color).
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("H3").Style = "Hyperlink"
$ws.Range("K3").Style = "Hyperlink"

$ws.Columns("K").Select()
